$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1464.125
$ws.Range("I15").Value = 1464.125
$ws.Range("K15").Value = 4392.375
$ws.Range("M15").Value = -4223.375
$ws.Range("H18").Value = 999.8333
$ws.Range("I18").Value = 499.66666
$ws.Range("K18").Value = 499.66666
$ws.Range("M18").Value = -215.66666
$ws.Range("H40").Value = 6852.0435
$ws.Range("I40").Value = 7750.375
$ws.Range("J40").Value = 6372.933
$ws.Range("K40").Value = 7750.375
$ws.Range("L40").Value = 6372.933
$ws.Range("M40").Value = -7575.375
$ws.Range("N40").Value = -6722.933
$ws.Range("H112").Value = 3608.4075
$ws.Range("J112").Value = 3840.7917
$ws.Range("L112").Value = 11522.3751
$ws.Range("N112").Value = -13738.3751
$ws.Range("H132").Value = 4560.357
$ws.Range("I132").Value = 4778.75
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 14336.25
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -11806.25
$ws.Range("N132").Value = -14810
$ws.Range("H138").Value = 6133.4907
$ws.Range("J138").Value = 7027.381
$ws.Range("L138").Value = 21082.143
$ws.Range("N138").Value = -31362.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2379.1538
$ws.Range("I61").Value = 1800.3478
$ws.Range("J61").Value = 6816.6665
$ws.Range("K61").Value = 1800.3478
$ws.Range("L61").Value = 6816.6665
$ws.Range("M61").Value = -1588.3478
$ws.Range("N61").Value = -7240.6665
$ws.Range("H136").Value = 2379.1538
$ws.Range("I136").Value = 1800.3478
$ws.Range("J136").Value = 6816.6665
$ws.Range("K136").Value = 5401.0434
$ws.Range("L136").Value = 20449.9995
$ws.Range("M136").Value = -2851.0434
$ws.Range("N136").Value = -25549.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 871.6
$ws.Range("I94").Value = 1034.75
$ws.Range("J94").Value = 685.1429000000001
$ws.Range("K94").Value = 1034.75
$ws.Range("L94").Value = 685.1429000000001
$ws.Range("M94").Value = -583.75
$ws.Range("N94").Value = -1587.1429
$ws.Range("H99").Value = 4203.6
$ws.Range("I99").Value = 3392.8333
$ws.Range("K99").Value = 3392.8333
$ws.Range("M99").Value = -1894.8333
$ws.Range("H122").Value = 73663.664
$ws.Range("J122").Value = 73663.664
$ws.Range("L122").Value = 73663.664
$ws.Range("N122").Value = -83463.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3502.5
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3502.5
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
$ws.Range("H122").Value = 2776.4
$ws.Range("I122").Value = 1527.15
$ws.Range("K122").Value = 4581.450000000001
$ws.Range("M122").Value = -2131.450000000001
$ws.Range("H139").Value = 74490
$ws.Range("J139").Value = 74490
$ws.Range("L139").Value = 74490
$ws.Range("N139").Value = -84770
$ws.Range("H141").Value = 103975.164
$ws.Range("J141").Value = 105516.06
$ws.Range("L141").Value = 105516.06
$ws.Range("N141").Value = -115876.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8500
$ws.Range("I56").Value = 8500
$ws.Range("K56").Value = 8500
$ws.Range("M56").Value = -7970
$ws.Range("H68").Value = 3334363.5
$ws.Range("J68").Value = 2501370
$ws.Range("L68").Value = 7504110
$ws.Range("N68").Value = -7505732
$ws.Range("H71").Value = 3334363.5
$ws.Range("J71").Value = 2501370
$ws.Range("L71").Value = 22512330
$ws.Range("N71").Value = -22520442
$ws.Range("H105").Value = 7900
$ws.Range("J105").Value = 7900
$ws.Range("L105").Value = 23700
$ws.Range("N105").Value = -28942
$ws.Range("H131").Value = 11181388
$ws.Range("I131").Value = 25718720
$ws.Range("J131").Value = 64604.59
$ws.Range("K131").Value = 77156160
$ws.Range("L131").Value = 193813.77
$ws.Range("M131").Value = -77151120
$ws.Range("N131").Value = -203893.77
$ws.Range("H132").Value = 1103286.4
$ws.Range("I132").Value = 252402.25
$ws.Range("J132").Value = 1670542.5
$ws.Range("K132").Value = 2271620.25
$ws.Range("L132").Value = 15034882.5
$ws.Range("M132").Value = -2269090.25
$ws.Range("N132").Value = -15039942.5
$ws.Range("H139").Value = 4809.778
$ws.Range("I139").Value = 2964.6667
$ws.Range("K139").Value = 8894.000100000001
$ws.Range("M139").Value = -3754.000100000001
$ws.Range("H140").Value = 4773.5557
$ws.Range("I140").Value = 4410.3335
$ws.Range("K140").Value = 13231.0005
$ws.Range("M140").Value = -8051.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5900
$ws.Range("J6").Value = 5900
$ws.Range("L6").Value = 5900
$ws.Range("N6").Value = -6126
$ws.Range("H10").Value = 1705000
$ws.Range("I10").Value = 2036000
$ws.Range("K10").Value = 2036000
$ws.Range("M10").Value = -2035831
$ws.Range("H11").Value = 8001666.5
$ws.Range("I11").Value = 9000000
$ws.Range("J11").Value = 7502500
$ws.Range("K11").Value = 9000000
$ws.Range("L11").Value = 7502500
$ws.Range("M11").Value = -8999861
$ws.Range("N11").Value = -7502778
$ws.Range("H13").Value = 3917.4285
$ws.Range("J13").Value = 3980.5
$ws.Range("L13").Value = 3980.5
$ws.Range("N13").Value = -4258.5
$ws.Range("H14").Value = 79228280
$ws.Range("I14").Value = 79228280
$ws.Range("K14").Value = 79228280
$ws.Range("M14").Value = -79228112
$ws.Range("H16").Value = 5900
$ws.Range("J16").Value = 5900
$ws.Range("L16").Value = 5900
$ws.Range("N16").Value = -6400
$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H21").Value = 1000000000
$ws.Range("I21").Value = 1000000000
$ws.Range("K21").Value = 1000000000
$ws.Range("M21").Value = -999999827
$ws.Range("H22").Value = 3500
$ws.Range("I22").Value = 3500
$ws.Range("K22").Value = 3500
$ws.Range("M22").Value = -2971
$ws.Range("H24").Value = 34107.145
$ws.Range("I24").Value = 9545.454
$ws.Range("K24").Value = 9545.454
$ws.Range("M24").Value = -9372.454
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H30").Value = 1000000000
$ws.Range("I30").Value = 1000000000
$ws.Range("K30").Value = 1000000000
$ws.Range("M30").Value = -999999895
$ws.Range("H33").Value = 25701172
$ws.Range("I33").Value = 12345
$ws.Range("K33").Value = 12345
$ws.Range("M33").Value = -12093
$ws.Range("H36").Value = 4049.75
$ws.Range("I36").Value = 2399
$ws.Range("K36").Value = 2399
$ws.Range("M36").Value = -1914
$ws.Range("H40").Value = 47212.25
$ws.Range("I40").Value = 44924.5
$ws.Range("K40").Value = 44924.5
$ws.Range("M40").Value = -44773.5
$ws.Range("H102").Value = 2463.6128
$ws.Range("I102").Value = 1640.3914
$ws.Range("K102").Value = 1640.3914
$ws.Range("M102").Value = -18.39139999999998
$ws.Range("H113").Value = 441919.78
$ws.Range("I113").Value = 627337.8
$ws.Range("K113").Value = 627337.8
$ws.Range("M113").Value = -625167.8
$ws.Range("H126").Value = 4143.4287
$ws.Range("I126").Value = 2077
$ws.Range("K126").Value = 6231
$ws.Range("M126").Value = -3761

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 630990.7
$ws.Range("I7").Value = 8943
$ws.Range("J7").Value = 913739.6
$ws.Range("K7").Value = 8943
$ws.Range("L7").Value = 913739.6
$ws.Range("M7").Value = -8831
$ws.Range("N7").Value = -913963.6
$ws.Range("H16").Value = 512.61536
$ws.Range("J16").Value = 600.6667
$ws.Range("L16").Value = 600.6667
$ws.Range("N16").Value = -940.6667
$ws.Range("H126").Value = 630990.7
$ws.Range("I126").Value = 8943
$ws.Range("J126").Value = 913739.6
$ws.Range("K126").Value = 26829
$ws.Range("L126").Value = 2741218.8
$ws.Range("M126").Value = -24359
$ws.Range("N126").Value = -2746158.8
$ws.Range("H132").Value = 5378.125
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 45000
$ws.Range("J99").Value = 45000
$ws.Range("L99").Value = 45000
$ws.Range("N99").Value = -50990
$ws.Range("H122").Value = 26319658
$ws.Range("I122").Value = 33336728
$ws.Range("J122").Value = 5643.25
$ws.Range("K122").Value = 100010184
$ws.Range("L122").Value = 16929.75
$ws.Range("M122").Value = -100007734
$ws.Range("N122").Value = -21829.75
$ws.Range("H132").Value = 32801.793
$ws.Range("I132").Value = 2846.52
$ws.Range("K132").Value = 8539.559999999999
$ws.Range("M132").Value = -6009.559999999999
